# Updates the cryptos list sheet1 with refreshed price/volume figures
# (and corrects the Bittensor/Filecoin/RenderToken row ordering).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.159.89"
$ws.Range("E2").Value = "  -1.51%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.483.43"
$ws.Range("E3").Value = "  -1.28%  "
# Row 4
$ws.Range("E4").Value = "  -0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.37"
$ws.Range("E5").Value = "  -2.77%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.78"
$ws.Range("E6").Value = "  -1.70%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "
# Row 8
$ws.Range("E8").Value = "  -1.64%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.482.95"
$ws.Range("E9").Value = "  -1.65%  "
# Row 10
$ws.Range("E10").Value = "  -3.49%  "
# Row 11
$ws.Range("E11").Value = "  -0.83%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  -2.27%  "
# Row 13
$ws.Range("E13").Value = "  -2.82%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.909.75"
$ws.Range("E14").Value = "  -1.93%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.121.09"
$ws.Range("E15").Value = "  -1.42%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.96"
$ws.Range("E16").Value = "  -3.89%  "
# Row 17
$ws.Range("E17").Value = "  -2.56%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.476.16"
$ws.Range("E18").Value = "  -1.92%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.61"
$ws.Range("E19").Value = "  -4.33%  "
# Row 20
$ws.Range("E20").Value = "  -2.20%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.92"
$ws.Range("E21").Value = "  -1.72%  "
# Row 22
$ws.Range("E22").Value = "  +0.09%  "
# Row 23
$ws.Range("E23").Value = "  -3.70%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.77"
# Row 25
$ws.Range("E25").Value = "  -2.80%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.12%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("E27").Value = "  -1.54%  "
# Row 28
$ws.Range("E28").Value = "  -2.90%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.39"
$ws.Range("E29").Value = "  +0.78%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0747"
$ws.Range("E30").Value = "  -2.58%  "
# Row 31
$ws.Range("E31").Value = "  -2.09%  "
# Row 32
$ws.Range("E32").Value = "  -3.40%  "
# Row 33
$ws.Range("E33").Value = "  +4.42%  "
# Row 35
$ws.Range("E35").Value = "  -0.02%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.05"
$ws.Range("E36").Value = "  -1.88%  "
# Row 37
$ws.Range("E37").Value = "  -4.15%  "
# Row 38
$ws.Range("E38").Value = "  -2.11%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.63"
$ws.Range("E39").Value = "  -0.86%  "
# Row 40
$ws.Range("E40").Value = "  -3.82%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.796"
$ws.Range("E41").Value = "  -1.16%  "
# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "274.48"
$ws.Range("E42").Value = "  -2.97%  "
# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.44"
$ws.Range("E43").Value = "  -4.34%  "
# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("E44").Value = "  -0.16%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.594"
$ws.Range("E45").Value = "  -2.01%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.99"
$ws.Range("E46").Value = "  -4.83%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0908"
$ws.Range("E47").Value = "  -1.72%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0489"
# Row 49
$ws.Range("E49").Value = "  -2.55%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.98"
$ws.Range("E50").Value = "  -2.20%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.733.85"
$ws.Range("E51").Value = "  -1.48%  "
